$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived NATMI metrics for Rps19-C5ar1 (rows 2-17).
# Maps row number -> hashtable of column letter -> new numeric value.
$rowUpdates = @{
    2 = @{ "G" = 92.32855733333334; "H" = 276.985672; "I" = 0.1878287475723421; "J" = 0.1878287475723421; "M" = 0.1137376666666667; "N" = 0.341213; "O" = 0.003048021899328029; "P" = 0.003048021899328029; "Q" = 10.50123467779289; "R" = 94.511112100136; "S" = 0.0005725061359238553; "T" = 0.0005725061359238553 }
    3 = @{ "G" = 92.32855733333334; "H" = 276.985672; "I" = 0.1878287475723421; "J" = 0.1878287475723421; "O" = 0.0144044366216848; "P" = 0.0144044366216848; "Q" = 49.62706120945334; "R" = 446.64355088508; "S" = 0.002705567290136236; "T" = 0.002705567290136236 }
    4 = @{ "G" = 92.32855733333334; "H" = 276.985672; "I" = 0.1878287475723421; "J" = 0.1878287475723421; "K" = 1; "L" = 0.3333333333333333; "M" = 0.05518366666666667; "N" = 0.165551; "O" = 0.001478850669393178; "P" = 0.001478850669393178; "Q" = 5.095028331696889; "R" = 45.855254985272; "S" = 0.0002777706690786405; "T" = 0.0002777706690786405 }
    5 = @{ "G" = 92.32855733333334; "H" = 276.985672; "I" = 0.1878287475723421; "J" = 0.1878287475723421; "M" = 36.608813; "N" = 109.826439; "O" = 0.981068690809594; "P" = 0.9810686908095939; "Q" = 3380.038889975779; "R" = 30420.35000978201; "S" = 0.1842729034772034; "T" = 0.1842729034772034 }
    6 = @{ "I" = 0.2518455369783797; "J" = 0.2518455369783797; "M" = 0.1137376666666667; "N" = 0.341213; "O" = 0.003048021899328029; "P" = 0.003048021899328029; "Q" = 14.08032114650678; "R" = 126.722890318561; "S" = 0.0007676307119581285; "T" = 0.0007676307119581285 }
    7 = @{ "I" = 0.2518455369783797; "J" = 0.2518455369783797; "O" = 0.0144044366216848; "P" = 0.0144044366216848; "S" = 0.003627693075859248; "T" = 0.003627693075859248 }
    8 = @{ "I" = 0.2518455369783797; "J" = 0.2518455369783797; "K" = 1; "L" = 0.3333333333333333; "M" = 0.05518366666666667; "N" = 0.165551; "O" = 0.001478850669393178; "P" = 0.001478850669393178; "Q" = 6.831542895860778; "R" = 61.483886062747; "S" = 0.0003724419409441613; "T" = 0.0003724419409441613 }
    9 = @{ "I" = 0.2518455369783797; "J" = 0.2518455369783797; "M" = 36.608813; "N" = 109.826439; "O" = 0.981068690809594; "P" = 0.9810686908095939; "Q" = 4532.041661651921; "R" = 40788.37495486729; "S" = 0.2470777712496182; "T" = 0.2470777712496182 }
    10 = @{ "G" = 174.4662783333333; "H" = 523.3988350000001; "I" = 0.3549257510290025; "J" = 0.3549257510290025; "M" = 0.1137376666666667; "N" = 0.341213; "O" = 0.003048021899328029; "P" = 0.003048021899328029; "Q" = 19.84338740965056; "R" = 178.590486686855; "S" = 0.001081821461771847; "T" = 0.001081821461771847 }
    11 = @{ "G" = 174.4662783333333; "H" = 523.3988350000001; "I" = 0.3549257510290025; "J" = 0.3549257510290025; "O" = 0.0144044366216848; "P" = 0.0144044366216848; "Q" = 93.77649693555834; "R" = 843.9884724200251; "S" = 0.005112505486101147; "T" = 0.005112505486101147 }
    12 = @{ "G" = 174.4662783333333; "H" = 523.3988350000001; "I" = 0.3549257510290025; "J" = 0.3549257510290025; "K" = 1; "L" = 0.3333333333333333; "M" = 0.05518366666666667; "N" = 0.165551; "O" = 0.001478850669393178; "P" = 0.001478850669393178; "Q" = 9.627688948120555; "R" = 86.64920053308501; "S" = 0.000524882184494117; "T" = 0.000524882184494117 }
    13 = @{ "G" = 174.4662783333333; "H" = 523.3988350000001; "I" = 0.3549257510290025; "J" = 0.3549257510290025; "M" = 36.608813; "N" = 109.826439; "O" = 0.981068690809594; "P" = 0.9810686908095939; "Q" = 6387.003358310953; "R" = 57483.03022479857; "S" = 0.3482065418966354; "T" = 0.3482065418966354 }
    14 = @{ "G" = 100.965814; "H" = 302.897442; "I" = 0.2053999644202756; "J" = 0.2053999644202756; "M" = 0.1137376666666667; "N" = 0.341213; "O" = 0.003048021899328029; "P" = 0.003048021899328029; "Q" = 11.48361609746066; "R" = 103.352544877146; "S" = 0.0006260635896741981; "T" = 0.0006260635896741981 }
    15 = @{ "G" = 100.965814; "H" = 302.897442; "I" = 0.2053999644202756; "J" = 0.2053999644202756; "O" = 0.0144044366216848; "P" = 0.0144044366216848; "Q" = 54.26962985406999; "R" = 488.4266686866299; "S" = 0.002958670769588174; "T" = 0.002958670769588174 }
    16 = @{ "G" = 100.965814; "H" = 302.897442; "I" = 0.2053999644202756; "J" = 0.2053999644202756; "K" = 1; "L" = 0.3333333333333333; "M" = 0.05518366666666667; "N" = 0.165551; "O" = 0.001478850669393178; "P" = 0.001478850669393178; "Q" = 5.571663824504665; "R" = 50.14497442054199; "S" = 0.0003037558748762596; "T" = 0.0003037558748762596 }
    17 = @{ "G" = 100.965814; "H" = 302.897442; "I" = 0.2053999644202756; "J" = 0.2053999644202756; "M" = 36.608813; "N" = 109.826439; "O" = 0.981068690809594; "P" = 0.9810686908095939; "Q" = 3696.238604118782; "R" = 33266.14743706903; "S" = 0.201511474186137; "T" = 0.201511474186137 }
}

foreach ($row in $rowUpdates.Keys) {
    $cols = $rowUpdates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
